# Adds the new "city" / "company_name" columns to the PF1.0 sheet, matching
# the credit_simulation_data.xlsx test-fixture update: reading the ID-card
# photos (anverso/reverso) also pulled through a couple of extra profile
# fields that now get captured/asserted in the spreadsheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PF1.0")

# --- Header row (row 1): two new trailing columns ---------------------
$ws.Range("S1").Value = "city"
$ws.Range("T1").Value = "company_name"

# --- Data row (row 2): values for the new columns ----------------------
# "city" mirrors the existing "commune" value (K2) for this sample row.
$ws.Range("S2").Value = $ws.Range("K2").Value()
$ws.Range("T2").Value = "N/A"

# --- Column formatting: column T (20) gets an explicit custom width ----
$ws.Columns.Item(20).ColumnWidth = 16

# --- View state: scroll so column J is first, select U8 ----------------
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollColumn = 10
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("U8").Select() | Out-Null
